# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
# Updates rows 119-121 on the "India Super League" sheet: the three most
# recent fixtures shift up one slot and a brand new fixture is appended,
# matching the upstream data refresh described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 119 (was id 7749470 / Kerala Blasters vs East Bengal Club) ---
# becomes the former row 120 data (Chennaiyin FC vs Jamshedpur FC) with
# refreshed closing-line odds.
$ws.Range("B119").Value = 7751760
$ws.Range("E119").Value = 45386.45833333334
$ws.Range("F119").Value = "Chennaiyin FC"
$ws.Range("G119").Value = "Jamshedpur FC"
$ws.Range("K119").Value = 2.3
$ws.Range("L119").Value = 3.5
$ws.Range("M119").Value = 2.9
$ws.Range("N119").Value = 2.3
$ws.Range("O119").Value = 3.5
$ws.Range("P119").Value = 2.9
$ws.Range("Q119").Value = -0.25
$ws.Range("R119").Value = 2
$ws.Range("S119").Value = 1.8
$ws.Range("T119").Value = 2.75
$ws.Range("U119").Value = 1.925
$ws.Range("V119").Value = 1.875

# --- Row 120 (was id 7751760 / Chennaiyin FC vs Jamshedpur FC) ---
# becomes the former row 121 data (FC Goa vs Hyderabad FC) with refreshed
# closing-line odds.
$ws.Range("B120").Value = 7751761
$ws.Range("E120").Value = 45387.45833333334
$ws.Range("F120").Value = "FC Goa"
$ws.Range("G120").Value = "Hyderabad FC"
$ws.Range("K120").Value = 1.142
$ws.Range("L120").Value = 7.5
$ws.Range("M120").Value = 19
$ws.Range("N120").Value = 1.125
$ws.Range("O120").Value = 8
$ws.Range("P120").Value = 21
$ws.Range("Q120").Value = -2.25
$ws.Range("R120").Value = 1.95
$ws.Range("S120").Value = 1.85
$ws.Range("T120").Value = 3.25
$ws.Range("U120").Value = 1.975
$ws.Range("V120").Value = 1.825

# --- Row 121 (was id 7751761 / FC Goa vs Hyderabad FC) ---
# becomes a brand-new fixture: Punjab FC vs Mohun Bagan SG.
$ws.Range("B121").Value = 7749762
$ws.Range("E121").Value = 45388.35416666666
$ws.Range("F121").Value = "Punjab FC"
$ws.Range("G121").Value = "Mohun Bagan SG"
$ws.Range("K121").Value = 4.333
$ws.Range("L121").Value = 3.6
$ws.Range("M121").Value = 1.75
$ws.Range("N121").Value = 4.333
$ws.Range("O121").Value = 3.6
$ws.Range("P121").Value = 1.75
$ws.Range("Q121").Value = 0.75
$ws.Range("R121").Value = 1.8
$ws.Range("S121").Value = 2
$ws.Range("T121").Value = 2.5
$ws.Range("U121").Value = 1.85
$ws.Range("V121").Value = 1.95
